$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the species-record data (everything except the
# shared location/admin columns, which are already identical across the
# three rows) among rows 16, 17 and 18:
#   new row16 <- old row17
#   new row17 <- old row18
#   new row18 <- old row16
# Only columns A,B,D,E,F,G,H,I,J,K,L,P,Q,R actually change value; the rest
# (C,T,U,V,W,Y,Z,AA,AB,AD,AE,AF,AG,AT,AW,AX,AY) are identical across the
# three rows already and are left untouched.

# Helper: write a literal (non-formula) text value into a cell without
# letting Excel auto-convert a numeric-looking string into a Number --
# build it as a formula returning the text, then bake the formula down to
# a plain value in place via Copy + PasteSpecial(values). This avoids the
# "quote prefix" / Text-number-format style that a direct
# Range.Value = "10" assignment would otherwise stamp onto the cell.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# ---------------- Row 16 (<- old row 17 content) ----------------
$ws.Range("A16").Value = 111837675
$ws.Range("B16").Value = 103288
Set-TextValue $ws.Range("D16") "LC"
$ws.Range("E16").Value = 221144
Set-TextValue $ws.Range("F16") "Grönpyrola"
Set-TextValue $ws.Range("G16") "Pyrola chlorantha"
Set-TextValue $ws.Range("H16") "Sw."
Set-TextValue $ws.Range("I16") "10"
Set-TextValue $ws.Range("J16") "plantor/tuvor"
$ws.Range("L16").Formula = '=""'
Set-TextValue $ws.Range("P16") "Brotorp, Långsjön, Sm"
$ws.Range("Q16").Value = 575781.9606960951
$ws.Range("R16").Value = 6404546.96767282

# ---------------- Row 17 (<- old row 18 content) ----------------
$ws.Range("A17").Value = 111837741
$ws.Range("B17").Value = 90658
Set-TextValue $ws.Range("D17") "NT"
$ws.Range("E17").Value = 4361
Set-TextValue $ws.Range("F17") "Orange taggsvamp"
Set-TextValue $ws.Range("G17") "Hydnellum aurantiacum"
Set-TextValue $ws.Range("H17") "(Batsch:Fr.) P.Karst."
Set-TextValue $ws.Range("I17") "15"
Set-TextValue $ws.Range("J17") "fruktkroppar"
$ws.Range("L17").ClearContents()
Set-TextValue $ws.Range("P17") "Brotorp, hyggeskant, Sm"
$ws.Range("Q17").Value = 575653.9215098171
$ws.Range("R17").Value = 6404506.688862759

# ---------------- Row 18 (<- old row 16 content) ----------------
$ws.Range("A18").Value = 111837705
$ws.Range("B18").Value = 90662
Set-TextValue $ws.Range("D18") "LC"
$ws.Range("E18").Value = 4363
Set-TextValue $ws.Range("F18") "Zontaggsvamp"
Set-TextValue $ws.Range("G18") "Hydnellum concrescens"
Set-TextValue $ws.Range("H18") "(Pers.) Banker"
Set-TextValue $ws.Range("I18") "10"
Set-TextValue $ws.Range("J18") "fruktkroppar"
Set-TextValue $ws.Range("P18") "Brotorp, Långsjön, Sm"
$ws.Range("Q18").Value = 575795.3141537429
$ws.Range("R18").Value = 6404518.948622406
